$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A to hold a date value for each movement row.
$ws.Range("A1:A2").EntireColumn.Insert()

# Fill the new column A with the movement date for both existing rows.
$ws.Range("A1").Value = "27/04/2018"
$ws.Range("A2").Value = "27/04/2018"
